$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("K12").Style = "Comma"
Write-Output "styled"
